# Saldo.xlsx update — apply the account-balance refresh described in the diff:
#  - ANA's balance drops to 166175.48, and a new GUILHERME (005637820) row for
#    100000 is added right after her
#  - POLYANNA's balance rises to 55656.99, and FLK (005883672) is removed
#  - RODRIGO (005152037) and REDRAU (008008723) are removed and replaced by
#    PAULO (004419141, 37199.38) and VICTOR (008032413, 20000)
#  - JORGE (005599726) is removed
#  - PEDRO (005880628), the old PAULO (004419141, 6508.04), DANIELA
#    (004452507), EMMANUELLE (004206790) and RODRIGO (004272426) rows are
#    removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Remove-AccountRow($account) {
    $hit = $ws.Columns.Item(1).Find($account)
    if ($hit) {
        $ws.Rows.Item($hit.Row()).Delete()
    }
}

function Set-AccountBalance($account, $newValue) {
    $hit = $ws.Columns.Item(1).Find($account)
    $hit.Offset(0, 2).Value = $newValue
}

function Insert-AccountRowAfter($afterAccount, $account, $name, $value) {
    $hit = $ws.Columns.Item(1).Find($afterAccount)
    $newRow = $hit.Row() + 1
    $ws.Rows.Item($newRow).Insert()
    $ws.Cells.Item($newRow, 1).Value = "'" + $account
    $ws.Cells.Item($newRow, 2).Value = $name
    $ws.Cells.Item($newRow, 3).Value = $value
}

# --- Rows removed outright ---
Remove-AccountRow "005883672"   # FLK
Remove-AccountRow "005152037"   # RODRIGO (30754.71)
Remove-AccountRow "008008723"   # REDRAU
Remove-AccountRow "005599726"   # JORGE
Remove-AccountRow "005880628"   # PEDRO (8435.7)
Remove-AccountRow "004419141"   # PAULO (6508.04, the old one)
Remove-AccountRow "004452507"   # DANIELA
Remove-AccountRow "004206790"   # EMMANUELLE
Remove-AccountRow "004272426"   # RODRIGO (1612.69)

# --- Balance updates on existing rows ---
Set-AccountBalance "004432579" 166175.48   # ANA
Set-AccountBalance "004389994" 55656.99    # POLYANNA

# --- New rows ---
Insert-AccountRowAfter "004432579" "005637820" "GUILHERME" 100000
Insert-AccountRowAfter "004322719" "004419141" "PAULO" 37199.38
Insert-AccountRowAfter "004419141" "008032413" "VICTOR" 20000
